$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.551.08"
$ws.Range("E2").Value = "  -0.36%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.627.16"
$ws.Range("E3").Value = "  +0.38%  "

# Row 4
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.60"
$ws.Range("E5").Value = "  -0.44%  "

# Row 6
$ws.Range("E6").Value = "  -0.92%  "

# Row 7
$ws.Range("E7").Value = "  -0.19%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.247"
$ws.Range("E8").Value = "  +0.21%  "

# Row 9
$ws.Range("E9").Value = "  -0.18%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.26"
$ws.Range("E10").Value = "  -0.17%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0856"
$ws.Range("E11").Value = "  -0.40%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.853.88"
$ws.Range("E12").Value = "  +0.26%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.639.45"
$ws.Range("E13").Value = "  +1.15%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.06"
$ws.Range("E14").Value = "  +0.15%  "

# Row 15
$ws.Range("E15").Value = "  +0.40%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.88"
$ws.Range("E16").Value = "  -1.15%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.566.95"
$ws.Range("E17").Value = "  -0.32%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "234.49"
$ws.Range("E18").Value = "  +0.59%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.80"
$ws.Range("E19").Value = "  +1.16%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0726"
$ws.Range("E20").Value = "  -0.55%  "

# Row 21
$ws.Range("E21").Value = "  -0.16%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.33"
$ws.Range("E22").Value = "  -1.41%  "

# Row 23
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.18"
$ws.Range("E23").Value = "  +0.95%  "

# Row 24
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.20"
$ws.Range("E24").Value = "  -1.59%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.09"
$ws.Range("E25").Value = "  +0.19%  "

# Row 26
$ws.Range("E26").Value = "  -0.21%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.08"
$ws.Range("E27").Value = "  +0.48%  "

# Row 28
$ws.Range("E28").Value = "  -0.47%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.69"
$ws.Range("E29").Value = "  +0.36%  "

# Row 30
$ws.Range("E30").Value = "  -0.54%  "

# Row 31
$ws.Range("E31").Value = "  -0.58%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.526.27"
$ws.Range("E32").Value = "  +4.69%  "

# Row 33
$ws.Range("E33").Value = "  +0.48%  "

# Row 34
$ws.Range("E34").Value = "  +0.69%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.54"
$ws.Range("E35").Value = "  +3.44%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.41"
$ws.Range("E36").Value = "  -0.51%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.569"
$ws.Range("E37").Value = "  -0.04%  "

# Row 38
$ws.Range("E38").Value = "  -0.74%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.836"
$ws.Range("E39").Value = "  -0.23%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.86"
$ws.Range("E40").Value = "  -0.80%  "

# Row 41
$ws.Range("E41").Value = "  -0.21%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.21"
$ws.Range("E42").Value = "  +0.47%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.766.77"
$ws.Range("E43").Value = "  +0.36%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.01"
$ws.Range("E44").Value = "  +1.94%  "

# Row 45
$ws.Range("E45").Value = "  -0.72%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.908"
$ws.Range("E46").Value = "  -4.31%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.78"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.52"
$ws.Range("E48").Value = "  +1.14%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0101"
$ws.Range("E49").Value = "  -3.86%  "

# Row 50
$ws.Range("E50").Value = "  -0.40%  "

# Row 51
$ws.Range("E51").Value = "  +0.16%  "
